$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('E2').Value = "2026-02-07 06:17:34"
$ws.Range('H2').Value = "'91%"
$ws.Range('O2').Value = "-1.7 °C"
$ws.Range('E3').Value = "2026-02-07 06:17:37"
$ws.Range('H3').Value = "'93%"
$ws.Range('K3').Value = "-0.1 MJ/m2"
$ws.Range('N3').Value = "-7.2 °C 5:39 TU"
$ws.Range('O3').Value = "-5.9 °C"
$ws.Range('E4').Value = "2026-02-07 06:17:39"
$ws.Range('J4').Value = "1001.3 hPa"
$ws.Range('N4').Value = "9.8 °C 5:48 TU"
$ws.Range('O4').Value = "11.3 °C"
$ws.Range('E5').Value = "2026-02-07 06:17:42"
$ws.Range('H5').Value = "'74%"
$ws.Range('J5').Value = "1001.3 hPa"
$ws.Range('L5').Value = "15.1 km/h - 202º 5:47 TU"
$ws.Range('N5').Value = "5.8 °C 5:53 TU"
$ws.Range('O5').Value = "8.6 °C"
$ws.Range('E6').Value = "2026-02-07 06:17:44"
$ws.Range('H6').Value = "'58%"
$ws.Range('J6').Value = "1002.9 hPa"
$ws.Range('E7').Value = "2026-02-07 06:17:47"
$ws.Range('J7').Value = "1002.7 hPa"
$ws.Range('E8').Value = "2026-02-07 06:17:49"
$ws.Range('H8').Value = "'95%"
$ws.Range('O8').Value = "3.9 °C"
$ws.Range('E9').Value = "2026-02-07 06:17:51"
$ws.Range('O9').Value = "1.5 °C"
$ws.Range('E10').Value = "2026-02-07 06:17:54"
$ws.Range('O10').Value = "7.3 °C"
$ws.Range('E11').Value = "2026-02-07 06:17:56"
$ws.Range('J11').Value = "1005.6 hPa"
$ws.Range('E12').Value = "2026-02-07 06:17:58"
$ws.Range('H12').Value = "'71%"
$ws.Range('E13').Value = "2026-02-07 06:18:01"
$ws.Range('H13').Value = "'87%"
$ws.Range('M13').Value = "11.1 °C 5:46 TU"
$ws.Range('O13').Value = "7.7 °C"
$ws.Range('E14').Value = "2026-02-07 06:18:03"
$ws.Range('H14').Value = "'73%"
$ws.Range('N14').Value = "-6.8 °C 5:40 TU"
$ws.Range('O14').Value = "-5.7 °C"
$ws.Range('E15').Value = "2026-02-07 06:18:06"
$ws.Range('H15').Value = "'87%"
$ws.Range('J15').Value = "1001.6 hPa"
$ws.Range('N15').Value = "2.9 °C 5:59 TU"
$ws.Range('O15').Value = "6.1 °C"
$ws.Range('E16').Value = "2026-02-07 06:18:08"
$ws.Range('O16').Value = "2.9 °C"
$ws.Range('E17').Value = "2026-02-07 06:18:11"
$ws.Range('J17').Value = "1004.9 hPa"
$ws.Range('N17').Value = "2.4 °C 5:34 TU"
$ws.Range('E18').Value = "2026-02-07 06:18:13"
$ws.Range('N18').Value = "-9.6 °C 5:58 TU"
$ws.Range('O18').Value = "-7.5 °C"
$ws.Range('E19').Value = "2026-02-07 06:18:16"
$ws.Range('J19').Value = "1006.2 hPa"
$ws.Range('N19').Value = "1.8 °C 5:57 TU"
$ws.Range('O19').Value = "4.3 °C"
$ws.Range('E20').Value = "2026-02-07 06:18:18"
$ws.Range('H20').Value = "'83%"
$ws.Range('K20').Value = "-0.1 MJ/m2"
$ws.Range('N20').Value = "-5.7 °C 5:58 TU"
$ws.Range('O20').Value = "-4.7 °C"
$ws.Range('E21').Value = "2026-02-07 06:18:21"
$ws.Range('H21').Value = "'76%"
$ws.Range('J21').Value = "1001.9 hPa"
$ws.Range('N21').Value = "2.1 °C 5:32 TU"
$ws.Range('O21').Value = "6.4 °C"
$ws.Range('E22').Value = "2026-02-07 06:18:23"
$ws.Range('H22').Value = "'91%"
$ws.Range('K22').Value = "-0.1 MJ/m2"
$ws.Range('M22').Value = "8.3 °C 5:58 TU"
$ws.Range('O22').Value = "6.1 °C"
$ws.Range('E23').Value = "2026-02-07 06:18:26"
$ws.Range('N23').Value = "6.9 °C 5:57 TU"
$ws.Range('E24').Value = "2026-02-07 06:18:28"
$ws.Range('H24').Value = "'81%"
$ws.Range('K24').Value = "-0.1 MJ/m2"
$ws.Range('L24').Value = "64.4 km/h - 335º 5:47 TU"
$ws.Range('O24').Value = "10.3 °C"
$ws.Range('E25').Value = "2026-02-07 06:18:31"
$ws.Range('J25').Value = "1005.2 hPa"
$ws.Range('E26').Value = "2026-02-07 06:18:33"
$ws.Range('N26').Value = "-5.4 °C 5:53 TU"
$ws.Range('O26').Value = "-2.1 °C"
$ws.Range('E27').Value = "2026-02-07 06:18:36"
$ws.Range('H27').Value = "'94%"
$ws.Range('J27').Value = "1001.2 hPa"
$ws.Range('L27').Value = "17.3 km/h - 12º 5:41 TU"
$ws.Range('M27').Value = "11.3 °C 5:47 TU"
$ws.Range('O27').Value = "8.8 °C"
$ws.Range('E28').Value = "2026-02-07 06:18:38"
$ws.Range('H28').Value = "'89%"
$ws.Range('J28').Value = "1004.1 hPa"
$ws.Range('O28').Value = "3.1 °C"
$ws.Range('E29').Value = "2026-02-07 06:18:40"
$ws.Range('L29').Value = "38.9 km/h - 272º 5:38 TU"
$ws.Range('E30').Value = "2026-02-07 06:18:43"
$ws.Range('H30').Value = "'82%"
$ws.Range('E31').Value = "2026-02-07 06:18:45"
$ws.Range('N31').Value = "3.0 °C 5:30 TU"
$ws.Range('E32').Value = "2026-02-07 06:18:48"
$ws.Range('J32').Value = "1004.4 hPa"
$ws.Range('E33').Value = "2026-02-07 06:18:50"
$ws.Range('H33').Value = "'91%"
$ws.Range('N33').Value = "4.6 °C 5:41 TU"
$ws.Range('O33').Value = "6.8 °C"
$ws.Range('E34').Value = "2026-02-07 06:18:52"
$ws.Range('E35').Value = "2026-02-07 06:18:55"
$ws.Range('O35').Value = "-5.9 °C"
$ws.Range('E36').Value = "2026-02-07 06:18:57"
$ws.Range('J36').Value = "1006.7 hPa"
$ws.Range('N36').Value = "4.0 °C 5:37 TU"
